$wb = $excel.ActiveWorkbook

# Add the new "LOINC" sheet. Worksheets.Add() with no args inserts the new
# sheet immediately before the active sheet, i.e. as the very first tab --
# exactly where it needs to land per the target layout.
$loinc = $wb.Worksheets.Add()
$loinc.Name = "LOINC"

# Borrow the header formatting (bold font + fill) already used by the other
# "CODE"/"<col2 header>" header cells elsewhere in the workbook (e.g. the
# Tier sheet's A1/B1), so the new header row matches styles s=1 / s=2.
$tier = $wb.Worksheets.Item("Tier")
$tier.Range("A1").Copy()
$loinc.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$tier.Range("B1").Copy()
$loinc.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$loinc.Range("C1").PasteSpecial(-4122)   # xlPasteFormats (same style as B1)

# Header row
$loinc.Range("A1").Value = "CODE"
$loinc.Range("B1").Value = "DESCRIPTION"
$loinc.Range("C1").Value = "HL7"

# Column A (codes) filled first, then column B (descriptions), then column C
# (formulas) -- matches the shared-string insertion order of the target.
$loinc.Range("A2").Value = "TPG_PAND"
$loinc.Range("A3").Value = "TPG_POP_GRP"
$loinc.Range("A4").Value = "TPG_TIER"

$loinc.Range("B2").Value = "Priority Group - Pandemic"
$loinc.Range("B3").Value = "Priority Group - Population Group"
$loinc.Range("B4").Value = "Priority Group - Tier"

$loinc.Range("C2").Formula = '=A2&"^"&SUBSTITUTE(B2, "&", "\T\")&"^99TPG"'
$loinc.Range("C3").Formula = '=A3&"^"&SUBSTITUTE(B3, "&", "\T\")&"^99TPG"'
$loinc.Range("C4").Formula = '=A4&"^"&SUBSTITUTE(B4, "&", "\T\")&"^99TPG"'

# Column widths approximating the authored sheet (autofit-style widths).
$loinc.Range("A1").ColumnWidth = 20.88
$loinc.Range("B1").ColumnWidth = 30.45
$loinc.Range("C1").ColumnWidth = 13.74

# Leave the cursor where the author last left it on this sheet.
$loinc.Range("C6").Select() | Out-Null
